# Update absenteeism data rows 2-11 with new values as per the target dataset.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 50183
$ws.Range("B2").Value = "Ana Vitória Azevedo"
$ws.Range("C2").Value = "Engenharia"
$ws.Range("D2").Value = "Consulta médica"
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 45099
$ws.Range("G2").Value = 7111.08

# Row 3
$ws.Range("A3").Value = 79933
$ws.Range("B3").Value = "Maria Alice da Conceição"
$ws.Range("C3").Value = "Marketing"
$ws.Range("D3").Value = "Consulta médica"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 45090
$ws.Range("G3").Value = 3631

# Row 4
$ws.Range("A4").Value = 57422
$ws.Range("B4").Value = "Dr. André Ferreira"
$ws.Range("C4").Value = "Vendas"
$ws.Range("D4").Value = "Viagem de negócios"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 45090
$ws.Range("G4").Value = 7957.54

# Row 5
$ws.Range("A5").Value = 56965
$ws.Range("B5").Value = "João Miguel da Paz"
$ws.Range("C5").Value = "Financeiro"
$ws.Range("D5").Value = "Viagem de negócios"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 45091
$ws.Range("G5").Value = 8633.549999999999

# Row 6
$ws.Range("A6").Value = 43896
$ws.Range("B6").Value = "Dr. João Felipe Peixoto"
$ws.Range("C6").Value = "Recursos Humanos"
$ws.Range("D6").Value = "Viagem de negócios"
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 45082
$ws.Range("G6").Value = 5157.24

# Row 7
$ws.Range("A7").Value = 92561
$ws.Range("B7").Value = "Alana Cardoso"
$ws.Range("C7").Value = "P&D"
$ws.Range("D7").Value = "Viagem de negócios"
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 45090
$ws.Range("G7").Value = 7439.26

# Row 8
$ws.Range("A8").Value = 71231
$ws.Range("B8").Value = "Marcela da Costa"
$ws.Range("C8").Value = "P&D"
$ws.Range("D8").Value = "Viagem de negócios"
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = 45103
$ws.Range("G8").Value = 9292.51

# Row 9
$ws.Range("A9").Value = 22904
$ws.Range("B9").Value = "Arthur Pires"
$ws.Range("C9").Value = "Vendas"
$ws.Range("D9").Value = "Viagem de negócios"
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = 45081
$ws.Range("G9").Value = 3783.1

# Row 10
$ws.Range("A10").Value = 7717
$ws.Range("B10").Value = "Sra. Stephany Barros"
$ws.Range("C10").Value = "Vendas"
$ws.Range("D10").Value = "Doença"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 45078
$ws.Range("G10").Value = 10602.61

# Row 11
$ws.Range("A11").Value = 12013
$ws.Range("B11").Value = "Bruna Correia"
$ws.Range("C11").Value = "Atendimento ao Cliente"
$ws.Range("D11").Value = "Doença"
$ws.Range("E11").Value = 6
$ws.Range("F11").Value = 45094
$ws.Range("G11").Value = 7753.07
